$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D (shifts existing D:K data to F:M)
# across the full used-row range so existing column widths / formats
# defined on entire columns follow correctly.
$ws.Range("D5:E102").EntireColumn.Insert()

# Copy number/date formatting from the (now-shifted) old column D, which
# landed in column F, into the two freshly inserted blank columns D:E so
# the new cells pick up the same date / thousands-style formatting as
# the rest of each table, restricted to each tables row-range so blank
# header rows are left untouched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns with the latest two quarters of data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 406400
$ws.Range("E8").Value = 409900
$ws.Range("D9").Value = 379300
$ws.Range("E9").Value = 385300
$ws.Range("D10").Value = 27100
$ws.Range("E10").Value = 24600
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 403900
$ws.Range("E17").Value = 409300
$ws.Range("D18").Value = 2500
$ws.Range("E18").Value = 600
$ws.Range("D20").Value = -37300
$ws.Range("E20").Value = 6100
$ws.Range("D21").Value = -16800
$ws.Range("E21").Value = 24800
$ws.Range("D22").Value = 4200
$ws.Range("E22").Value = 4400
$ws.Range("D23").Value = -39000
$ws.Range("E23").Value = 2300
$ws.Range("D24").Value = 600
$ws.Range("E24").Value = 900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -39600
$ws.Range("E26").Value = 1400
$ws.Range("D27").Value = -39600
$ws.Range("E27").Value = 1400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29:E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 37300
$ws.Range("E32").Value = -6100
$ws.Range("D33").Value = -39600
$ws.Range("E33").Value = 1400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -39600
$ws.Range("E35").Value = 1400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 21600
$ws.Range("E41").Value = 27000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 163400
$ws.Range("E43").Value = 172100
$ws.Range("D44").Value = 296800
$ws.Range("E44").Value = 285300
$ws.Range("D45").Value = 9600
$ws.Range("E45").Value = 13000
$ws.Range("D46").Value = 491400
$ws.Range("E46").Value = 497400
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 674400
$ws.Range("E48").Value = 670200
$ws.Range("D49").Value = 17800
$ws.Range("E49").Value = 16800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 14000
$ws.Range("E52").Value = 22900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1197600
$ws.Range("E54").Value = 1207300
$ws.Range("D57").Value = 163600
$ws.Range("E57").Value = 141700
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 57200
$ws.Range("E59").Value = 54900
$ws.Range("D60").Value = 220800
$ws.Range("E60").Value = 196600
$ws.Range("D61").Value = 189100
$ws.Range("E61").Value = 218000
$ws.Range("D62").Value = 252500
$ws.Range("E62").Value = 218400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 662400
$ws.Range("E66").Value = 633000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -269200
$ws.Range("E72").Value = -229600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 535200
$ws.Range("E76").Value = 574300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -39600
$ws.Range("E81").Value = 1400
$ws.Range("D83").Value = 18000
$ws.Range("E83").Value = 18100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 46900
$ws.Range("E89").Value = 1700
$ws.Range("D91").Value = -22300
$ws.Range("E91").Value = -8700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -22300
$ws.Range("E94").Value = -8700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -30000
$ws.Range("E100").Value = -5000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -5400
$ws.Range("E102").Value = -12000
